$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Noten_ID" column inserted before the existing table (A1) and a new
# "Anzahl Klausur in Schuljahr" column appended after the existing table (N1).
$ws.Range("A1").Value = "Noten_ID"
$ws.Range("N1").Value = "Anzahl Klausur in Schuljahr"

# Give the new last column a wider, custom width (as the others already have).
$ws.Columns.Item(14).ColumnWidth = 22.39

# Scroll the view so column E is the left-most visible column, with M1 kept
# as the active cell (mirrors the author's on-screen state when saving).
$ws.Range("M1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
